# Fix list level numbering: paragraphs that were incorrectly bumped to the
# second outline level (OOXML a:pPr/@lvl="1", i.e. COM IndentLevel 2) when
# they should have stayed at the top level (lvl="0" / IndentLevel 1).
#
# We walk every slide, every shape's text (including grouped shapes and
# table cells), and for every paragraph currently at IndentLevel 2, demote
# it back to IndentLevel 1 - leaving any other level untouched.

function Fix-TextRangeLevels($textRange) {
    if ($textRange -eq $null) { return }
    $count = $textRange.Paragraphs().Count
    for ($i = 1; $i -le $count; $i++) {
        $para = $textRange.Paragraphs($i, 1)
        if ($para.IndentLevel -eq 2) {
            $para.IndentLevel = 1
        }
    }
}

function Fix-Shape($shape) {
    if ($shape.HasTextFrame -and $shape.TextFrame.HasText) {
        Fix-TextRangeLevels $shape.TextFrame.TextRange
    }
    if ($shape.HasTable) {
        $tbl = $shape.Table
        for ($r = 1; $r -le $tbl.Rows.Count; $r++) {
            for ($c = 1; $c -le $tbl.Columns.Count; $c++) {
                $cell = $tbl.Cell($r, $c)
                Fix-TextRangeLevels $cell.Shape.TextFrame.TextRange
            }
        }
    }
    if ($shape.Type -eq 6) {
        # msoGroup
        for ($gi = 1; $gi -le $shape.GroupItems.Count; $gi++) {
            Fix-Shape $shape.GroupItems.Item($gi)
        }
    }
}

$p = $ppt.ActivePresentation
for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $slide = $p.Slides.Item($si)
    for ($shi = 1; $shi -le $slide.Shapes.Count; $shi++) {
        Fix-Shape $slide.Shapes.Item($shi)
    }
}
